# update model2 comparison data and visualizations
# Row 8 ("covenext_large") comparison metrics were refreshed with new
# validation/testing numbers, and the validation_loss cell (B8) is
# highlighted with a new left-aligned, word-wrapped, white-filled style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refresh the comparison metrics for the "covenext_large" row (row 8) ---
$ws.Range("B8").Value = 0.125
$ws.Range("C8").Value = 0.96637
$ws.Range("D8").Value = 0.10765
$ws.Range("E8").Value = 0.96
# F8 (input_features) is unchanged.

# --- give B8 (validation_loss) a distinct highlight style ---
$b8 = $ws.Range("B8")
$b8.Font.Color = 0          # explicit black (RGB 0,0,0)
$b8.Font.Size = 11          # bump from the default 10pt
$b8.Interior.Color = 16777215   # solid white fill
$b8.HorizontalAlignment = -4131 # xlLeft
$b8.WrapText = $true
